# Fix XML validation issue: w:shd elements with a w:fill attribute but no
# w:val attribute. Word COM always emits a texture/pattern (w:val) alongside
# a fill color; the source docs were missing it. Setting Shading.Texture to
# wdTextureNone (0) on the already-shaded header cells forces the w:val
# attribute to be (re)written without touching the fill color itself.

$wdTextureNone = 0
$wdColorAutomatic = -16777216

$lightGray = 13882323   # RGB D3D3D3
$darkGray  = 11119017   # RGB A9A9A9

$d = $word.ActiveDocument

for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $table = $d.Tables.Item($t)
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        for ($c = 1; $c -le $table.Columns.Count; $c++) {
            $cell = $table.Cell($r, $c)
            $bg = $cell.Shading.BackgroundPatternColor
            if ($bg -eq $lightGray -or $bg -eq $darkGray) {
                $cell.Shading.Texture = $wdTextureNone
            }
        }
    }
}

Write-Output "done"
